$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "AUD"
$ws.Range("G9").Value = "CHF"
$ws.Range("G14").Value = "CHF"
$ws.Range("G17").Value = "AUD"
$ws.Range("G18").Value = "GBP"
$ws.Range("G19").Value = "AUD"
$ws.Range("G20").Value = "COP"
$ws.Range("G22").Value = "AUD"
$ws.Range("G26").Value = "COP"
$ws.Range("G30").Value = "AUD"
$ws.Range("G32").Value = "COP"
$ws.Range("G41").Value = "AUD"
$ws.Range("G47").Value = "COP"
$ws.Range("G50").Value = "AUD"
$ws.Range("G51").Value = "AUD"
$ws.Range("G53").Value = "AUD"
$ws.Range("G55").Value = "CHF"
$ws.Range("G57").Value = "AUD"
$ws.Range("G64").Value = "AUD"
$ws.Range("G65").Value = "AUD"
$ws.Range("G74").Value = "AUD"
$ws.Range("G75").Value = "COP"
$ws.Range("G86").Value = "AUD"
$ws.Range("G89").Value = "COP"
$ws.Range("G94").Value = "JPY"
$ws.Range("G101").Value = "AUD"
$ws.Range("G111").Value = "COP"
$ws.Range("G113").Value = "GBP"
$ws.Range("G116").Value = "COP"
$ws.Range("G120").Value = "AUD"
$ws.Range("G121").Value = "COP"
$ws.Range("G130").Value = "AUD"
$ws.Range("G132").Value = "AUD"
$ws.Range("G135").Value = "AUD"
$ws.Range("G136").Value = "AUD"
$ws.Range("G148").Value = "AUD"
$ws.Range("G154").Value = "AUD"
$ws.Range("G157").Value = "AUD"
$ws.Range("G161").Value = "AUD"
$ws.Range("G168").Value = "GBP"
$ws.Range("G173").Value = "AUD"
$ws.Range("G175").Value = "AUD"
$ws.Range("G177").Value = "AUD"
$ws.Range("G178").Value = "AUD"
$ws.Range("G181").Value = "COP"
$ws.Range("G182").Value = "COP"
$ws.Range("G187").Value = "AUD"
$ws.Range("G189").Value = "GBP"
$ws.Range("G192").Value = "AUD"
$ws.Range("G194").Value = "AUD"
$ws.Range("G195").Value = "GBP"
$ws.Range("G197").Value = "CHF"
$ws.Range("G201").Value = "CHF"
$ws.Range("G204").Value = "AUD"
$ws.Range("G217").Value = "AUD"
$ws.Range("G221").Value = "CHF"
$ws.Range("G225").Value = "COP"
$ws.Range("G237").Value = "AUD"
$ws.Range("G241").Value = "AUD"
$ws.Range("G252").Value = "COP"
$ws.Range("G255").Value = "AUD"
$ws.Range("G265").Value = "GBP"
$ws.Range("G266").Value = "AUD"
$ws.Range("G270").Value = "AUD"
$ws.Range("G271").Value = "AUD"
$ws.Range("G273").Value = "USD"
$ws.Range("I273").Value = "Buenos días, habla Sara Gómez de la Mesa de Divisas, es un buen momento para conversar. Estoy revisando opciones para comprar `$74,356 dolares, ¿qué condiciones tienen hoy? Esta es la mejor tasa que puedo ofrecer ahora mismo, `$3,820.62. Me gustaría que revisaras si puede mejorar la propuesta. Esta es la mejor tasa que puedo ofrecer ahora mismo, `$3,810.94. Si me ofrece una mejora mínima, podemos concretar la operación Con la tendencia actual del mercado, 3,801.5 es una tasa muy competitiva Si mejora la tasa en al menos 15 pesos, cerramos Le puedo ofrecer una tasa de 3,792.44 si realiza la operación en los próximos minutos Podría aceptar si ajusta la tasa a unos puntos a mi favor Podemos realizar la transacción a una tasa preferencial de 3.792.44 pesos, válida solo durante los próximos minutos De acuerdo, me parece bien la oferta, avancemos Gracias. Gracias. ¿Quieres hacer referencia? 646."
$ws.Range("G284").Value = "AUD"
$ws.Range("G286").Value = "COP"
$ws.Range("G289").Value = "AUD"
$ws.Range("G299").Value = "AUD"
$ws.Range("G311").Value = "COP"
$ws.Range("G317").Value = "AUD"
$ws.Range("G319").Value = "COP"
$ws.Range("G321").Value = "COP"
$ws.Range("G324").Value = "AUD"
$ws.Range("G331").Value = "AUD"
$ws.Range("G332").Value = "AUD"
$ws.Range("G336").Value = "AUD"
$ws.Range("G346").Value = "JPY"
$ws.Range("G348").Value = "GBP"
$ws.Range("G351").Value = "CHF"
$ws.Range("G356").Value = "AUD"
$ws.Range("G358").Value = "AUD"
$ws.Range("G359").Value = "COP"
$ws.Range("G364").Value = "AUD"
$ws.Range("G366").Value = "AUD"
$ws.Range("G367").Value = "AUD"
$ws.Range("G368").Value = "COP"
$ws.Range("G372").Value = "AUD"
$ws.Range("G373").Value = "AUD"
$ws.Range("G384").Value = "COP"
$ws.Range("G385").Value = "AUD"
$ws.Range("G387").Value = "JPY"
$ws.Range("G395").Value = "COP"
$ws.Range("G396").Value = "COP"
$ws.Range("G397").Value = "AUD"
$ws.Range("G400").Value = "AUD"
$ws.Range("G408").Value = "COP"
$ws.Range("G409").Value = "CHF"
$ws.Range("G413").Value = "AUD"
$ws.Range("G420").Value = "AUD"
$ws.Range("G426").Value = "COP"
$ws.Range("G428").Value = "AUD"
$ws.Range("G432").Value = "GBP"
$ws.Range("G436").Value = "AUD"
$ws.Range("G442").Value = "AUD"
$ws.Range("G447").Value = "AUD"
$ws.Range("G450").Value = "COP"
$ws.Range("G453").Value = "AUD"
$ws.Range("G458").Value = "AUD"
$ws.Range("G462").Value = "AUD"
$ws.Range("G466").Value = "CHF"
$ws.Range("G473").Value = "AUD"
$ws.Range("G485").Value = "AUD"
$ws.Range("G492").Value = "AUD"
$ws.Range("G495").Value = "AUD"
$ws.Range("G501").Value = "COP"
